# Restore C10 on the "Rules" sheet to its saved value of 1 (numeric).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1.0
